$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 32, shifting the existing rows 32:56 down to 33:57.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new weekly record.
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 44778
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100108
$ws.Range("H32").Value = "Tropicales y subtropicales"
$ws.Range("I32").Value = 100108003
$ws.Range("J32").Value = "Maracuyá"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 20
$ws.Range("N32").Value = 35000
$ws.Range("O32").Value = 35000
$ws.Range("P32").Value = 35000
$ws.Range("Q32").Value = "$/caja 18 kilos"
$ws.Range("R32").Value = "Región de Arica y Parinacota"
$ws.Range("S32").Value = 1944
$ws.Range("T32").Value = 18
